$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue 'D2' '302.49'
Set-TextValue 'E2' '-5.87%'
Set-TextValue 'D3' '35.15'
Set-TextValue 'E3' '-2.94%'
Set-TextValue 'D4' '5.041'
Set-TextValue 'E4' '-1.63%'
Set-TextValue 'D5' '0.07910'
Set-TextValue 'E5' '-2.84%'
Set-TextValue 'D6' '1.945'
Set-TextValue 'E6' '-9.59%'
Set-TextValue 'D7' '7.751'
Set-TextValue 'E7' '-3.60%'
Set-TextValue 'D8' '4.024'
Set-TextValue 'E8' '-2.85%'
Set-TextValue 'D9' '2.870'
Set-TextValue 'D10' '0.9231'
Set-TextValue 'E10' '-0.51%'
Set-TextValue 'D11' '0.1198'
Set-TextValue 'E11' '18.55%'
Set-TextValue 'D12' '0.1838'
Set-TextValue 'E12' '-2.51%'
Set-TextValue 'D13' '0.09400'
Set-TextValue 'E13' '2.64%'
Set-TextValue 'D14' '0.03544'
Set-TextValue 'E14' '-1.30%'
Set-TextValue 'D15' '0.09870'
Set-TextValue 'E15' '-0.54%'
Set-TextValue 'D16' '0.001396'
Set-TextValue 'E16' '-2.47%'
Set-TextValue 'D17' '0.005852'
Set-TextValue 'E17' '3.58%'
Set-TextValue 'D18' '3.492'
Set-TextValue 'E18' '1.19%'
Set-TextValue 'D19' '0.3443'
Set-TextValue 'E19' '2.07%'
Set-TextValue 'E20' '-0.16%'
Set-TextValue 'D21' '5.037'
Set-TextValue 'E21' '-0.49%'
Set-TextValue 'D22' '0.2401'
Set-TextValue 'E22' '9.61%'
Set-TextValue 'D23' '0.04508'
Set-TextValue 'E23' '-2.19%'
Set-TextValue 'E24' '-2.28%'
Set-TextValue 'D25' '0.004572'
Set-TextValue 'E25' '-3.34%'
Set-TextValue 'E26' '-3.90%'
Set-TextValue 'E27' '-6.91%'
Set-TextValue 'D39' '0.01895'
Set-TextValue 'E39' '-6.15%'
Set-TextValue 'D40' '0.04708'
Set-TextValue 'E40' '-5.58%'
Set-TextValue 'D41' '0.007632'
Set-TextValue 'E41' '-2.29%'
Set-TextValue 'D42' '0.009545'
Set-TextValue 'E42' '22.26%'
Set-TextValue 'E43' '-5.61%'
Set-TextValue 'E44' '1.36%'
Set-TextValue 'D45' '0.01120'
Set-TextValue 'E45' '-8.07%'
Set-TextValue 'D46' '0.00006023'
Set-TextValue 'E46' '-7.11%'
Set-TextValue 'E47' '-0.05%'
Set-TextValue 'E49' '-31.42%'
Set-TextValue 'D50' '0.00002101'
Set-TextValue 'E50' '-0.05%'
Set-TextValue 'D51' '0.0002001'
Set-TextValue 'E51' '-0.05%'
